# issue #5: stock data from json to db
#
# The "股票" (stock) worksheet gains three new columns:
#   - "category"    inserted right after "property_category" (becomes column I)
#   - "source_file" appended after "legislator_id"            (becomes column M)
#   - "index"       appended after "source_file"               (becomes column N)
#
# category    = "normal"     for every data row
# source_file = "tmp2e4a1"   for every data row
# index       = the same value already stored in column A for that row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 6

# --- Insert the new "category" column as column I (index 9) ---------------
# This shifts the existing I/J/K columns (date/legislator_name/legislator_id)
# one place to the right, to J/K/L, and inherits their cell formatting.
$ws.Columns.Item(9).Insert()

$ws.Cells.Item(1, 9).Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# --- Append "source_file" (M) and "index" (N) columns ----------------------
# Inserting columns immediately to the right of the current last used column
# (L) makes Excel copy the formatting of column L into the new columns.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmp2e4a1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
